$d = $word.ActiveDocument

# 1 & 7. Shorten title / bold restatement of the title
$d.Content.Find.Execute(
    "Play Beer Bonanza for Free - A Festive Slot Game with Exciting Features",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Beer Bonanza for Free", 2)

# 2. "Unique and engaging cluster payouts" -> "Festive and colorful Oktoberfest theme"
$d.Content.Find.Execute(
    "Unique and engaging cluster payouts",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Festive and colorful Oktoberfest theme", 2)

# 3. "Festive and colorful Oktoberfest design" -> "Unique cluster payouts"
$d.Content.Find.Execute(
    "Festive and colorful Oktoberfest design",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Unique cluster payouts", 2)

# 4. "Special symbols offer free spins and multipliers" -> "Special symbols for free spins and multipliers"
$d.Content.Find.Execute(
    "Special symbols offer free spins and multipliers",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Special symbols for free spins and multipliers", 2)

# 5. "No jackpot available" -> "No jackpot feature"
$d.Content.Find.Execute(
    "No jackpot available",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "No jackpot feature", 2)

# 6. "Betting range may not appeal to high rollers" -> "Limited betting range"
$d.Content.Find.Execute(
    "Betting range may not appeal to high rollers",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Limited betting range", 2)

# 8. Italic summary sentence
$d.Content.Find.Execute(
    "Enjoy the Oktoberfest-themed Beer Bonanza slot game with cluster payouts, special symbols, and free spins. Play for free and experience the festival vibe!",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Enjoy the festive Oktoberfest theme and exciting bonus features in the free slot game Beer Bonanza.", 2)
